$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "darsh3@gmail.com"
$ws.Range("B2").Value = "sanj3@gmail.com"
$ws.Range("B3").Value = "harshi3@gmail.com"

$ws.Range("B6").Select()
